$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Parametrize the student data row (row 11) with template placeholder tokens ---
$ws.Range("A11").Value = "{StudentIndex}"
$ws.Range("B11").Value = "{StudentName}"
$ws.Range("C11").Value = "{OKRAvg}"
$ws.Range("D11").Value = "{LPR}"
$ws.Range("E11").Value = "{CourseGrade}"
$ws.Range("F11").Value = "{SemesterGrade}"

# --- Update the selected / active cell shown in the sheet view ---
$ws.Range("I19").Select() | Out-Null

# --- Slightly narrow column A (closest attainable width to 3.5703125 chars) ---
$ws.Columns.Item(1).ColumnWidth = 2.72
